$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.381320476531982
$ws.Range("B1").Value = 3.323147535324097
$ws.Range("C1").Value = 3.039629697799683
$ws.Range("D1").Value = 1.744714021682739
$ws.Range("E1").Value = 1.006046891212463
